# Atualiza pontuacoes e resultados das competicoes
#
# 1) "Geral": atualiza a pontuacao (coluna B) de varias equipes.
# 2) "Mes - Janeiro": recalcula a classificacao (ordenacao decrescente por
#    pontuacao), reescrevendo as colunas A (equipe) e B (pontuacao) na nova
#    ordem.

$wb    = $excel.ActiveWorkbook
$geral = $wb.Worksheets.Item("Geral")
$jan   = $wb.Worksheets.Item("Mes - Janeiro")

# --- 1) Novas pontuacoes na aba "Geral" (linha -> novo valor de B) ---------
$novasPontuacoes = @{
    2  = 64.5
    3  = 63.56
    4  = 47.86
    5  = 66.37
    7  = 53.66
    8  = 71.96
    9  = 56.05
    10 = 61.56
    12 = 47.86
    13 = 55.66
    14 = 62.56
    17 = 59.8
    18 = 61.16
    20 = 54.9
    21 = 68.06
}

foreach ($linha in $novasPontuacoes.Keys) {
    $geral.Cells.Item($linha, 2).Value = $novasPontuacoes[$linha]
}

# --- 2) Nova ordem de classificacao na aba "Mes - Janeiro" ------------------
# Cada item e a linha (na aba "Geral") da equipe que deve ocupar a
# respectiva posicao (linha 2..21) na aba "Mes - Janeiro", apos a
# atualizacao das pontuacoes acima.
$novaOrdemLinhasGeral = @(6, 8, 21, 19, 5, 2, 3, 14, 10, 18, 17, 9, 13, 15, 20, 11, 7, 4, 12, 16)

$linhaDestino = 2
foreach ($linhaOrigem in $novaOrdemLinhasGeral) {
    $equipe     = $geral.Cells.Item($linhaOrigem, 1).Value2
    $pontuacao  = $geral.Cells.Item($linhaOrigem, 2).Value2

    $jan.Cells.Item($linhaDestino, 1).Value = $equipe
    $jan.Cells.Item($linhaDestino, 2).Value = $pontuacao

    $linhaDestino = $linhaDestino + 1
}
